$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.06992366666666666
$ws.Range("H2").Value = 0.209771
$ws.Range("M2").Value = 0.06921533333333334
$ws.Range("N2").Value = 0.207646
$ws.Range("O2").Value = 0.005513961132583326
$ws.Range("P2").Value = 0.005513961132583327
$ws.Range("Q2").Value = 0.004839789896222222
$ws.Range("R2").Value = 0.04355810906599999
$ws.Range("S2").Value = 0.005513961132583326
$ws.Range("T2").Value = 0.005513961132583327

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 0.06992366666666666
$ws.Range("H3").Value = 0.209771
$ws.Range("O3").Value = 0.6731443835632516
$ws.Range("P3").Value = 0.6731443835632517
$ws.Range("Q3").Value = 0.5908415579893332
$ws.Range("R3").Value = 5.317574021903999
$ws.Range("S3").Value = 0.6731443835632516
$ws.Range("T3").Value = 0.6731443835632517

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 0.06992366666666666
$ws.Range("H4").Value = 0.209771
$ws.Range("M4").Value = 4.033719
$ws.Range("N4").Value = 12.101157
$ws.Range("O4").Value = 0.321341655304165
$ws.Range("P4").Value = 0.321341655304165
$ws.Range("Q4").Value = 0.282052422783
$ws.Range("R4").Value = 2.538471805047
$ws.Range("S4").Value = 0.321341655304165
$ws.Range("T4").Value = 0.321341655304165
